$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price observation was recorded ahead of the existing history for
# "Vega Modelo de Temuco - Acelga". Insert a fresh row at row 326 so every
# existing record (rows 326-341) shifts down by one (326-327 ... 341-342),
# matching the new used range A1:R342.
$ws.Rows.Item(326).Insert()

# Fill in the newly inserted row with the new observation's data. The
# descriptive columns (market/region/category/etc.) are identical to every
# other row in this block.
$ws.Cells.Item(326, 1).Value() = 10
$ws.Cells.Item(326, 2).Value() = "Vega Modelo de Temuco"
$ws.Cells.Item(326, 3).Value() = "La Araucanía"
$ws.Cells.Item(326, 4).Value() = 44753
$ws.Cells.Item(326, 5).Value() = 9
$ws.Cells.Item(326, 6).Value() = 100112009
$ws.Cells.Item(326, 7).Value() = "Acelga"
$ws.Cells.Item(326, 8).Value() = "Sin especificar"
$ws.Cells.Item(326, 9).Value() = "Primera"
$ws.Cells.Item(326, 10).Value() = 80
$ws.Cells.Item(326, 11).Value() = 12000
$ws.Cells.Item(326, 12).Value() = 12000
$ws.Cells.Item(326, 13).Value() = 12000
$ws.Cells.Item(326, 14).Value() = "`$/docena de atados (12 kilos)"
$ws.Cells.Item(326, 15).Value() = "Provincia de Cautín"
$ws.Cells.Item(326, 16).Value() = 1000
$ws.Cells.Item(326, 17).Value() = 12
$ws.Cells.Item(326, 18).Value() = "Hortaliza"
